$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.962.17"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "1.846.07"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "309.46"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4778"
$ws.Range("E7").Value = "  +1.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3673"
$ws.Range("E8").Value = "  +1.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07225"
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9281"
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.71"
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07737"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").Value = "1.836.41"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.345"
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.444"
$ws.Range("E15").Value = "  +1.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.91"
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.014"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008664"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "26.996.86"
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.47"
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.066"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.63"
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.923"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.89"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.24"
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.005"
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.17"
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.958"
$ws.Range("E29").Value = "  +1.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08859"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E31").Value = "  +5.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.175"
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7398"
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.499"
$ws.Range("E34").Value = "  +1.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.757"
$ws.Range("E35").Value = "  -2.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.113"
$ws.Range("E36").Value = "  +3.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01958"
$ws.Range("E37").Value = "  +1.47%  "
$ws.Range("E38").Value = "  +2.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.980"
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5229"
$ws.Range("E40").Value = "  +3.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.997"
$ws.Range("E41").Value = "  +1.98%  "
$ws.Range("E42").Value = "  +1.00%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.73"
$ws.Range("E43").Value = "  +6.61%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.245"
$ws.Range("E44").Value = "  +1.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4744"
$ws.Range("E45").Value = "  +2.18%  "
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.97"
$ws.Range("E47").Value = "  +3.39%  "
$ws.Range("E48").Value = "  +1.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "65.80"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06067"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8896"
$ws.Range("E51").Value = "  +3.74%  "
